$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.22690486907959
$ws.Range("B1").Value = 2.307052612304688
$ws.Range("C1").Value = 3.390167474746704
$ws.Range("D1").Value = 2.128296375274658
$ws.Range("E1").Value = 1.33695387840271
